$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target final state of rows 8..17 (A, B-name, C, D, E):
#   8  : 6  line7  14 11 TRUE
#   9  : 7  line8  16 9  FALSE
#   10 : 8  extr1  5  12 TRUE
#   11 : 9  extr2  5  9  TRUE
#   12 : 10 extr3  10 11 FALSE
#   13 : 11 extr4  7  8  FALSE
#   14 : 12 extr5  9  11 TRUE
#   15 : 13 extr6  7  11 FALSE
#   16 : 14 extr7  5  7  FALSE
#   17 : 15 extr8  8  5  FALSE

# First touch column A for rows 8..17 with plain values so the used range
# covers A1:E17 before we copy formatting over it (avoids Excel creating
# spurious intermediate cell-style combinations).
for ($r = 8; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
}

# Copy the existing bordered/bold/centered style used by column A (rows 2-7)
# onto the newly touched rows 8-17 so they match the rest of the table.
$ws.Range("A2:A7").Copy()
$ws.Range("A8:A17").PasteSpecial(-4122)

# Now write the actual values for the new/updated rows.
$names = @("line7", "line8", "extr1", "extr2", "extr3", "extr4", "extr5", "extr6", "extr7", "extr8")
$aVals = @(6, 7, 8, 9, 10, 11, 12, 13, 14, 15)
$cVals = @(14, 16, 5, 5, 10, 7, 9, 7, 5, 8)
$dVals = @(11, 9, 12, 9, 11, 8, 11, 11, 7, 5)
$eVals = @($true, $false, $true, $true, $false, $false, $true, $false, $false, $false)

for ($i = 0; $i -lt 10; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = $aVals[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
}
